# Generate Report for Handback
#
# This script records the results of a localization handback: for each
# localized file (zh-cn and de-de), the status moves from "Ready for
# handoff" to "Handed back: in sync with en-US", and the "Latest Target
# File" / "Latest Handback File" / "Latest Handback DateTime" columns are
# populated with the handback xliff info.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$statusText = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet: the per-locale status columns (E = zh-cn, F = de-de)
# ---------------------------------------------------------------------
$wsOverview.Range("E2").Value = $statusText
$wsOverview.Range("F2").Value = $statusText
$wsOverview.Range("E3").Value = $statusText
$wsOverview.Range("F3").Value = $statusText

# Widen the now-longer status columns to fit the new text.
$wsOverview.Columns.Item(5).ColumnWidth = 29.166666666666668
$wsOverview.Columns.Item(6).ColumnWidth = 29.166666666666668

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZhCn.Range("C2").Value = $statusText
$wsZhCn.Range("C3").Value = $statusText

# Latest Target File (I) - hyperlinked to the handed-back source markdown
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/8fbaa06a4c6077b31f4b47282fbc5a9429f70e45/e2e/8a5aa4f2-4315-488b-a160-d830e17efa0d.md", "", "", "8a5aa4f2-4315-488b-a160-d830e17efa0d.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/8fbaa06a4c6077b31f4b47282fbc5a9429f70e45/e2e/f5fcffff-7c7b-4436-a1ba-0c7e540837c7.md", "", "", "f5fcffff-7c7b-4436-a1ba-0c7e540837c7.md")

# Latest Handback File (J) - the xliff consumed for the handback
$wsZhCn.Range("J2").Value = "8a5aa4f2-4315-488b-a160-d830e17efa0d.0b33261dadb17f7d55ac0bc1ec5903e6d82b2a49.zh-cn.xlf"
$wsZhCn.Range("J3").Value = "f5fcffff-7c7b-4436-a1ba-0c7e540837c7.68e013114717f7c768d6e533ab3df7cbcd3bce1a.zh-cn.xlf"

# Latest Handback DateTime (K)
$wsZhCn.Range("K2").Value = "2016-08-13 21:02:30"
$wsZhCn.Range("K3").Value = "2016-08-13 21:02:30"

# Widen columns that now hold the longer status text / file names.
$wsZhCn.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsZhCn.Columns.Item(9).ColumnWidth = 39.166666666666664
$wsZhCn.Columns.Item(10).ColumnWidth = 39.166666666666664

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDeDe.Range("C2").Value = $statusText
$wsDeDe.Range("C3").Value = $statusText

# Latest Target File (I) - hyperlinked to the handed-back source markdown
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/8fbaa06a4c6077b31f4b47282fbc5a9429f70e45/e2e/8a5aa4f2-4315-488b-a160-d830e17efa0d.md", "", "", "8a5aa4f2-4315-488b-a160-d830e17efa0d.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/8fbaa06a4c6077b31f4b47282fbc5a9429f70e45/e2e/f5fcffff-7c7b-4436-a1ba-0c7e540837c7.md", "", "", "f5fcffff-7c7b-4436-a1ba-0c7e540837c7.md")

# Latest Handback File (J) - the xliff consumed for the handback
$wsDeDe.Range("J2").Value = "8a5aa4f2-4315-488b-a160-d830e17efa0d.0b33261dadb17f7d55ac0bc1ec5903e6d82b2a49.de-de.xlf"
$wsDeDe.Range("J3").Value = "f5fcffff-7c7b-4436-a1ba-0c7e540837c7.68e013114717f7c768d6e533ab3df7cbcd3bce1a.de-de.xlf"

# Latest Handback DateTime (K)
$wsDeDe.Range("K2").Value = "2016-08-13 21:02:40"
$wsDeDe.Range("K3").Value = "2016-08-13 21:02:40"

# Widen columns that now hold the longer status text / file names.
$wsDeDe.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsDeDe.Columns.Item(9).ColumnWidth = 39.166666666666664
$wsDeDe.Columns.Item(10).ColumnWidth = 39.166666666666664

Write-Output "Handback report generated."
